# Adds a new "Czech" worksheet (a copy of "Belgium") and fills in the
# Czech-specific test data, matching the "Added and Updated Czech Test
# Data" commit.

$wb = $excel.ActiveWorkbook
$belgium = $wb.Worksheets.Item("Belgium")

# Duplicate the Belgium sheet (same layout/styles) to seed the new Czech
# sheet, placing it right after Belgium.
$belgium.Copy($null, $belgium) | Out-Null
$czech = $wb.Worksheets.Item($wb.Worksheets.Count)
$czech.Name = "Czech"

# Fill in the Czech-specific values (market name + user story id).
$czech.Range("B2").Value = "Czech Market"
$czech.Range("B4").Value = "NGC-3477/T1734"

# Belgium is no longer the active tab, so drop its single-cell selection
# back to covering the whole data range.
$belgium.Range("A1:D10").Select() | Out-Null

# Make the new Czech sheet the active tab, with its own selection.
$czech.Select() | Out-Null
$czech.Range("B7").Select() | Out-Null
